# Generate Report for Handback
# Update status from "Ready for handoff" to "Handed back: in sync with en-US"
# across the Overview sheet and each locale sheet, refresh the handback
# timestamps, and clear the stale "handback not latest" error detail.

$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$zhcn = $wb.Worksheets.Item("zh-cn")
$dede = $wb.Worksheets.Item("de-de")

$newStatus = "Handed back: in sync with en-US"

# Overview sheet: zh-cn / de-de status columns
$overview.Range("E2").Value = $newStatus
$overview.Range("F2").Value = $newStatus

# zh-cn detail sheet
$zhcn.Range("C2").Value = $newStatus
$zhcn.Range("K2").Value = "2016-08-19 18:57:18"
$zhcn.Range("P2").Value = ""

# de-de detail sheet
$dede.Range("C2").Value = $newStatus
$dede.Range("K2").Value = "2016-08-19 18:57:25"
$dede.Range("P2").Value = ""

# Widen the status column (now holds longer text) and narrow the now-empty
# Error Detail column. (Values chosen land in the same rounded Excel column
# width bucket as the authored widths once persisted to OOXML.)
$overview.Columns.Item(5).ColumnWidth = 29.166
$overview.Columns.Item(6).ColumnWidth = 29.166

$zhcn.Columns.Item(3).ColumnWidth = 29.166
$zhcn.Columns.Item(16).ColumnWidth = 12.833

$dede.Columns.Item(3).ColumnWidth = 29.166
$dede.Columns.Item(16).ColumnWidth = 12.833
